# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2025-04-22 (serial 45769) to 2025-04-23 (serial 45770).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C43").Value = 45770
